$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; everything from old column B onward
# shifts one column to the right (B->C, C->D, ...).
$ws.Columns("B:B").Insert()

# Copy the header formatting (yellow fill header style) from column A's
# header cell onto the newly inserted header cell.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the new header its label and make the column as wide as column A.
$ws.Range("B1").Value = "比賽年份 Year of Competition"
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Leave the selection parked on A2, matching the saved workbook state.
[void]$ws.Range("A2").Select()
